$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5.737324714660645
$ws.Range("B1").Value = 4.805707931518555
$ws.Range("C1").Value = 3.380240201950073
$ws.Range("D1").Value = 2.008003234863281
$ws.Range("E1").Value = 1.629738450050354
